$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.802227087239764
$ws.Range("C2").Value = 6.049165449935222
$ws.Range("D2").Value = 5.974228524331156
$ws.Range("E2").Value = 16.49574201292873
$ws.Range("G2").Value = 25.2328075183318
$ws.Range("H2").Value = 13.38193240298047
$ws.Range("I2").Value = 19.05111207048686
$ws.Range("K2").Value = 8.795251567628915
$ws.Range("O2").Value = 19.90110053432614

$ws.Range("B3").Value = 8.41569684850084
$ws.Range("C3").Value = 5.811563294767859
$ws.Range("D3").Value = 5.85326309145764
$ws.Range("E3").Value = 15.56119451767119
$ws.Range("G3").Value = 25.27054076881181
$ws.Range("H3").Value = 13.43226909423748
$ws.Range("I3").Value = 19.1512246502861
$ws.Range("K3").Value = 8.444828403321209
$ws.Range("O3").Value = 19.97604942073842

$ws.Range("B4").Value = 8.169768105361106
$ws.Range("C4").Value = 5.659536688372878
$ws.Range("D4").Value = 5.779459808479268
$ws.Range("E4").Value = 14.96247296116223
$ws.Range("G4").Value = 25.30407212810599
$ws.Range("H4").Value = 13.46564885292902
$ws.Range("I4").Value = 19.21696406550554
$ws.Range("K4").Value = 8.220392489884865
$ws.Range("O4").Value = 20.02720517676504

$ws.Range("B5").Value = 8.067524334455729
$ws.Range("C5").Value = 5.596097088644188
$ws.Range("D5").Value = 5.749549694073675
$ws.Range("E5").Value = 14.71249588218556
$ws.Range("G5").Value = 25.32032753401625
$ws.Range("H5").Value = 13.47987250213483
$ws.Range("I5").Value = 19.24482554044844
$ws.Range("K5").Value = 8.12668101285931
$ws.Range("O5").Value = 20.04933838992627

$ws.Range("B6").Value = 8.050428828148764
$ws.Range("C6").Value = 5.585474928313584
$ws.Range("D6").Value = 5.744594516102799
$ws.Range("E6").Value = 14.67063432727695
$ws.Range("G6").Value = 25.32318275845808
$ws.Range("H6").Value = 13.48227182285498
$ws.Range("I6").Value = 19.24951662640783
$ws.Range("K6").Value = 8.110986834928772
$ws.Range("O6").Value = 20.05309116962562

$ws.Range("B7").Value = 8.168397218893958
$ws.Range("C7").Value = 5.658687063355626
$ws.Range("D7").Value = 5.779055698108491
$ws.Range("E7").Value = 14.95912555401343
$ws.Range("G7").Value = 25.30428088394173
$ws.Range("H7").Value = 13.46583816404335
$ws.Range("I7").Value = 19.21733547615293
$ws.Range("K7").Value = 8.21913766853786
$ws.Range("O7").Value = 20.02749846927237

$ws.Range("B8").Value = 8.670811889880671
$ws.Range("C8").Value = 5.968545447829428
$ws.Range("D8").Value = 5.932450171339288
$ws.Range("E8").Value = 16.17882526495157
$ws.Range("G8").Value = 25.24365883404134
$ws.Range("H8").Value = 13.39877484656722
$ws.Range("I8").Value = 19.08474371920385
$ws.Range("K8").Value = 8.676395411522133
$ws.Range("O8").Value = 19.92587386534885

$ws.Range("B9").Value = 9.582747070962501
$ws.Range("C9").Value = 6.525377183626143
$ws.Range("D9").Value = 6.234956180760538
$ws.Range("E9").Value = 18.44878882261968
$ws.Range("G9").Value = 25.2075681352119
$ws.Range("H9").Value = 13.28691539664609
$ws.Range("I9").Value = 18.85868252018835
$ws.Range("K9").Value = 9.496490200204592
$ws.Range("O9").Value = 19.76755997998172

$ws.Range("B10").Value = 10.20237979702595
$ws.Range("C10").Value = 6.901153105512184
$ws.Range("D10").Value = 6.455616513539707
$ws.Range("E10").Value = 20.08626243581957
$ws.Range("G10").Value = 25.23214713535756
$ws.Range("H10").Value = 13.21675455227679
$ws.Range("I10").Value = 18.71339201286911
$ws.Range("K10").Value = 10.04897030203451
$ws.Range("O10").Value = 19.67650850054713

$ws.Range("B11").Value = 10.47240948644817
$ws.Range("C11").Value = 7.064480104634161
$ws.Range("D11").Value = 6.555133428207874
$ws.Range("E11").Value = 20.78895828912996
$ws.Range("G11").Value = 25.2545088911846
$ws.Range("H11").Value = 13.18745636195474
$ws.Range("I11").Value = 18.65183416738782
$ws.Range("K11").Value = 10.28889575368012
$ws.Range("O11").Value = 19.64063020463406

$ws.Range("B12").Value = 10.57289826332373
$ws.Range("C12").Value = 7.12520668109584
$ws.Range("D12").Value = 6.592653273659427
$ws.Range("E12").Value = 21.0490186739538
$ws.Range("G12").Value = 25.2645875970926
$ws.Range("H12").Value = 13.17673923573123
$ws.Range("I12").Value = 18.62917814155828
$ws.Range("K12").Value = 10.37807331133244
$ws.Range("O12").Value = 19.62784541385059

$ws.Range("B13").Value = 10.55133562097661
$ws.Range("C13").Value = 7.112178401270578
$ws.Range("D13").Value = 6.584580671967304
$ws.Range("E13").Value = 20.99327791653261
$ws.Range("G13").Value = 25.26234530533871
$ws.Range("H13").Value = 13.17903055767658
$ws.Range("I13").Value = 18.63402836996155
$ws.Range("K13").Value = 10.3589424501268
$ws.Range("O13").Value = 19.63056312593069

$ws.Range("B14").Value = 10.48071246592506
$ws.Range("C14").Value = 7.069498737911071
$ws.Range("D14").Value = 6.558223742786656
$ws.Range("E14").Value = 20.81047438200144
$ws.Range("G14").Value = 25.25530578122132
$ws.Range("H14").Value = 13.18656708715787
$ws.Range("I14").Value = 18.64995710950613
$ws.Range("K14").Value = 10.29626624908822
$ws.Range("O14").Value = 19.63956230195838

$ws.Range("B15").Value = 10.43722212947939
$ws.Range("C15").Value = 7.043209395814102
$ws.Range("D15").Value = 6.542056691521115
$ws.Range("E15").Value = 20.69771693888543
$ws.Range("G15").Value = 25.25120368958693
$ws.Range("H15").Value = 13.19123261413548
$ws.Range("I15").Value = 18.65979924212688
$ws.Range("K15").Value = 10.25765585595188
$ws.Range("O15").Value = 19.64517907727918

$ws.Range("B16").Value = 10.18448881590138
$ws.Range("C16").Value = 6.890323680228767
$ws.Range("D16").Value = 6.449091932178947
$ws.Range("E16").Value = 20.03949327269582
$ws.Range("G16").Value = 25.23091104141346
$ws.Range("H16").Value = 13.21872202647991
$ws.Range("I16").Value = 18.71750643195947
$ws.Range("K16").Value = 10.03305784763529
$ws.Range("O16").Value = 19.67896514480667

$ws.Range("B17").Value = 10.02636380904209
$ws.Range("C17").Value = 6.794562397078844
$ws.Range("D17").Value = 6.391810213060949
$ws.Range("E17").Value = 19.62490965048337
$ws.Range("G17").Value = 25.22132874015671
$ws.Range("H17").Value = 13.23625718174107
$ws.Range("I17").Value = 18.75407125954367
$ws.Range("K17").Value = 9.892325583050637
$ws.Range("O17").Value = 19.7011144985423

$ws.Range("B18").Value = 9.934302673743263
$ws.Range("C18").Value = 6.738767559176917
$ws.Range("D18").Value = 6.358784401757159
$ws.Range("E18").Value = 19.38247818502258
$ws.Range("G18").Value = 25.21686960539656
$ws.Range("H18").Value = 13.246589331566
$ws.Range("I18").Value = 18.77552913950205
$ws.Range("K18").Value = 9.810308873021063
$ws.Range("O18").Value = 19.71437554216642

$ws.Range("B19").Value = 9.90294346980045
$ws.Range("C19").Value = 6.719754385194932
$ws.Range("D19").Value = 6.347590169664698
$ws.Range("E19").Value = 19.2997115235946
$ws.Range("G19").Value = 25.21554042815179
$ws.Range("H19").Value = 13.2501299092604
$ws.Range("I19").Value = 18.78286763356558
$ws.Range("K19").Value = 9.782356684974104
$ws.Range("O19").Value = 19.7189548954654

$ws.Range("B20").Value = 10.04331203863983
$ws.Range("C20").Value = 6.804830605628641
$ws.Range("D20").Value = 6.397916418386061
$ws.Range("E20").Value = 19.6694538327711
$ws.Range("G20").Value = 25.22223984808381
$ws.Range("H20").Value = 13.23436502696164
$ws.Range("I20").Value = 18.75013468624319
$ws.Range("K20").Value = 9.907417933120847
$ws.Range("O20").Value = 19.69870267268812

$ws.Range("B21").Value = 10.50150456435183
$ws.Range("C21").Value = 7.082065427914684
$ws.Range("D21").Value = 6.565970206262532
$ws.Range("E21").Value = 20.86433174618707
$ws.Range("G21").Value = 25.25732973186924
$ws.Range("H21").Value = 13.18434317393515
$ws.Range("I21").Value = 18.64526066765847
$ws.Range("K21").Value = 10.3147215524504
$ws.Range("O21").Value = 19.63689723581145

$ws.Range("B22").Value = 10.79064640075132
$ws.Range("C22").Value = 7.256706034839555
$ws.Range("D22").Value = 6.674823709261621
$ws.Range("E22").Value = 21.61009038022227
$ws.Range("G22").Value = 25.28965157213545
$ws.Range("H22").Value = 13.15385158671686
$ws.Range("I22").Value = 18.58053607320872
$ws.Range("K22").Value = 10.5711297408188
$ws.Range("O22").Value = 19.60117784483431

$ws.Range("B23").Value = 10.63728773647104
$ws.Range("C23").Value = 7.164104020628505
$ws.Range("D23").Value = 6.61682908742711
$ws.Range("E23").Value = 21.21527158482367
$ws.Range("G23").Value = 25.27154139030797
$ws.Range("H23").Value = 13.16992385174271
$ws.Range("I23").Value = 18.61473077538648
$ws.Range("K23").Value = 10.43518642049684
$ws.Range("O23").Value = 19.61981278284345

$ws.Range("B24").Value = 10.03565333042006
$ws.Range("C24").Value = 6.800190652715535
$ws.Range("D24").Value = 6.395156091957295
$ws.Range("E24").Value = 19.64932812497496
$ws.Range("G24").Value = 25.22182466612052
$ws.Range("H24").Value = 13.23521968845751
$ws.Range("I24").Value = 18.75191305218705
$ws.Range("K24").Value = 9.900598129525587
$ws.Range("O24").Value = 19.69979141760765

$ws.Range("B25").Value = 9.344534985065694
$ws.Range("C25").Value = 6.380444207906885
$ws.Range("D25").Value = 6.153215842418496
$ws.Range("E25").Value = 17.80813436839876
$ws.Range("G25").Value = 25.20839509366269
$ws.Range("H25").Value = 13.31506818483527
$ws.Range("I25").Value = 18.91619341162558
$ws.Range("K25").Value = 9.283211298095612
$ws.Range("O25").Value = 19.80597186381521
